$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-05-10 Saturday"; new = "2025-05-11 Sunday"},
    @{old = "803×3="; new = "721×6="},
    @{old = "129×8="; new = "628×5="},
    @{old = "510×5="; new = "592×4="},
    @{old = "738×6="; new = "898×5="},
    @{old = "256×3="; new = "851×5="},
    @{old = "605×5="; new = "721×4="},
    @{old = "152×3="; new = "153×9="},
    @{old = "390×7="; new = "607×4="},
    @{old = "857×6="; new = "535×9="},
    @{old = "483×5="; new = "659×4="},
    @{old = "334×6="; new = "738×5="},
    @{old = "285×4="; new = "485×2="},
    @{old = "645×7="; new = "999×9="},
    @{old = "103×9="; new = "845×3="},
    @{old = "938×8="; new = "499×8="},
    @{old = "254×5="; new = "177×9="},
    @{old = "663×5="; new = "427×8="},
    @{old = "551×7="; new = "977×3="},
    @{old = "214×5="; new = "780×5="},
    @{old = "804×3="; new = "551×2="},
    @{old = "940×9="; new = "323×2="},
    @{old = "273×4="; new = "458×6="},
    @{old = "404×9="; new = "446×9="},
    @{old = "413×6="; new = "869×4="},
    @{old = "523×9="; new = "701×3="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
